$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'292.08"
$ws.Range("E2").Value = "'-5.60%"
$ws.Range("D3").Value = "'40.09"
$ws.Range("E3").Value = "'-2.44%"
$ws.Range("E4").Value = "'-2.73%"
$ws.Range("D5").Value = "'0.07360"
$ws.Range("E5").Value = "'-3.98%"
$ws.Range("D6").Value = "'4.293"
$ws.Range("D7").Value = "'1.555"
$ws.Range("E7").Value = "'-9.08%"
$ws.Range("D8").Value = "'0.9183"
$ws.Range("E8").Value = "'0.37%"
$ws.Range("D9").Value = "'0.1187"
$ws.Range("E9").Value = "'-4.18%"
$ws.Range("D10").Value = "'0.1735"
$ws.Range("E10").Value = "'-4.20%"
$ws.Range("D11").Value = "'0.08743"
$ws.Range("E11").Value = "'-4.33%"
$ws.Range("D12").Value = "'0.04203"
$ws.Range("E12").Value = "'0.24%"
$ws.Range("D13").Value = "'0.1051"
$ws.Range("E13").Value = "'-0.08%"
$ws.Range("D14").Value = "'0.001280"
$ws.Range("E14").Value = "'-1.84%"
$ws.Range("D15").Value = "'0.005776"
$ws.Range("E15").Value = "'0.45%"
$ws.Range("D16").Value = "'3.406"
$ws.Range("E16").Value = "'1.81%"
$ws.Range("E18").Value = "'-0.70%"
$ws.Range("D19").Value = "'7.569"
$ws.Range("E19").Value = "'2.28%"
$ws.Range("E20").Value = "'-0.32%"
$ws.Range("D21").Value = "'0.2878"
$ws.Range("E21").Value = "'5.47%"
$ws.Range("D22").Value = "'0.03839"
$ws.Range("E22").Value = "'-4.57%"
$ws.Range("D23").Value = "'0.001281"
$ws.Range("E23").Value = "'0.89%"
$ws.Range("D24").Value = "'0.003890"
$ws.Range("E24").Value = "'-4.55%"
$ws.Range("D25").Value = "'0.0001282"
$ws.Range("E25").Value = "'-1.63%"
$ws.Range("D38").Value = "'0.02329"
$ws.Range("E38").Value = "'-7.77%"
$ws.Range("D39").Value = "'0.05027"
$ws.Range("E39").Value = "'-5.08%"
$ws.Range("D40").Value = "'0.007691"
$ws.Range("E40").Value = "'-1.94%"
$ws.Range("D41").Value = "'0.005110"
$ws.Range("E41").Value = "'172.49%"
$ws.Range("D42").Value = "'0.1272"
$ws.Range("E42").Value = "'-2.88%"
$ws.Range("D43").Value = "'0.007384"
$ws.Range("E43").Value = "'10.89%"
$ws.Range("D44").Value = "'0.007689"
$ws.Range("E44").Value = "'-5.54%"
$ws.Range("D45").Value = "'0.3168"
$ws.Range("E45").Value = "'3.22%"
$ws.Range("D46").Value = "'0.00006517"
$ws.Range("E46").Value = "'-4.25%"
$ws.Range("E47").Value = "'-0.21%"
$ws.Range("E48").Value = "'12.20%"
$ws.Range("E49").Value = "'35.58%"
$ws.Range("E50").Value = "'-0.21%"
$ws.Range("E51").Value = "'-0.21%"
